$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat "@" forces text (prevents Excel from auto-parsing values like
# "1.007" as numbers), and ClearFormats() afterwards restores the cell to the
# workbook default style (these cells carry no explicit style in the source).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "25.891.41"
Set-TextValue "E2" "  -1.73%  "
Set-TextValue "D3" "1.634.19"
Set-TextValue "E3" "  -1.96%  "
Set-TextValue "D4" "1.007"
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "215.46"
Set-TextValue "E5" "  -1.44%  "
Set-TextValue "D6" "0.5033"
Set-TextValue "E6" "  -2.49%  "
Set-TextValue "D7" "1.008"
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "0.2571"
Set-TextValue "E8" "  +0.23%  "
Set-TextValue "D9" "0.06430"
Set-TextValue "E9" "  -0.27%  "
Set-TextValue "D10" "19.52"
Set-TextValue "E10" "  -2.29%  "
Set-TextValue "D11" "0.07744"
Set-TextValue "E11" "  +1.17%  "
Set-TextValue "D12" "1.652.69"
Set-TextValue "E12" "  -0.64%  "
Set-TextValue "D13" "4.251"
Set-TextValue "E13" "  -1.76%  "
Set-TextValue "D14" "1.857.70"
Set-TextValue "E14" "  -1.99%  "
Set-TextValue "D15" "0.5452"
Set-TextValue "E15" "  -1.53%  "
Set-TextValue "D16" "0.0₅7939"
Set-TextValue "E16" "  -1.35%  "
Set-TextValue "D17" "63.51"
Set-TextValue "E17" "  -1.63%  "
Set-TextValue "D18" "25.868.77"
Set-TextValue "D19" "1.009"
Set-TextValue "E19" "  +0.10%  "
Set-TextValue "D20" "204.04"
Set-TextValue "E20" "  -3.04%  "
Set-TextValue "D21" "4.304"
Set-TextValue "E21" "  -2.21%  "
Set-TextValue "D22" "9.987"
Set-TextValue "E22" "  -1.29%  "
Set-TextValue "D23" "5.938"
Set-TextValue "E23" "  +0.70%  "
Set-TextValue "D24" "1.009"
Set-TextValue "E24" "  -0.02%  "
Set-TextValue "D25" "1.937"
Set-TextValue "E25" "  +10.88%  "
Set-TextValue "D26" "141.46"
Set-TextValue "E26" "  -2.26%  "
Set-TextValue "D27" "0.1154"
Set-TextValue "E27" "  -0.94%  "
Set-TextValue "D28" "15.77"
Set-TextValue "E28" "  -0.07%  "
Set-TextValue "D29" "6.766"
Set-TextValue "E29" "  -3.28%  "
Set-TextValue "D30" "0.05074"
Set-TextValue "E30" "  -3.74%  "
Set-TextValue "E31" "  -1.87%  "
Set-TextValue "D32" "3.263"
Set-TextValue "E32" "  -3.29%  "
Set-TextValue "D33" "3.192"
Set-TextValue "E33" "  -0.82%  "
Set-TextValue "D34" "1.544"
Set-TextValue "E34" "  -1.52%  "
Set-TextValue "D35" "2.341"
Set-TextValue "E35" "  -1.47%  "
Set-TextValue "D36" "0.8943"
Set-TextValue "E36" "  -3.64%  "
Set-TextValue "D37" "2.608"
Set-TextValue "D38" "0.5645"
Set-TextValue "E38" "  -1.33%  "
Set-TextValue "D39" "1.139.74"
Set-TextValue "E39" "  -0.93%  "
Set-TextValue "D40" "0.01558"
Set-TextValue "E40" "  -2.52%  "
Set-TextValue "D41" "2.566"
Set-TextValue "E41" "  -0.58%  "
Set-TextValue "D42" "1.009"
Set-TextValue "E42" "  +0.04%  "
Set-TextValue "D43" "5.635"
Set-TextValue "E43" "  -0.42%  "
Set-TextValue "E44" "  -3.31%  "
Set-TextValue "D45" "99.39"
Set-TextValue "E45" "  -0.62%  "
Set-TextValue "D46" "1.768.69"
Set-TextValue "D47" "0.0₈110"
Set-TextValue "E47" "  -2.37%  "
Set-TextValue "D48" "0.4517"
Set-TextValue "E48" "  +0.28%  "
Set-TextValue "D49" "1.011"
Set-TextValue "E49" "  +0.23%  "
Set-TextValue "D50" "54.79"
Set-TextValue "E50" "  -2.29%  "
Set-TextValue "D51" "0.05019"
Set-TextValue "E51" "  -1.72%  "
